# Apply the "Add files via upload" change to the "Особые события" sheet:
# - Row 32 becomes the last row of its group: restyle to the thin-bottom-border
#   variant (styles 8/10) used elsewhere in the workbook for group-end rows.
# - Ten new rows (33-42) of translated dialogue lines are appended.
# - The view selection/scroll position is moved to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Row 32: re-style as a group-terminator row (border styles 8/10). ---
# Sheet 1 already has this exact pattern (A/B = style 8, C/D/E = style 10) on
# row 9, so copy its formatting across rather than re-deriving it by hand.
$ws1.Range("A9:E9").Copy()
$ws2.Range("A32:E32").PasteSpecial(-4122)

# --- Row 33: new script group "SCRIPT/D25P11A/um2304.ssb" ---
$ws2.Range("A33").Value = "SCRIPT/D25P11A/um2304.ssb"
$ws2.Range("B33").Value = 250
$ws2.Range("C33").Value = " Wrooooah![K] We failed to get\nthrough the dungeon! It is a travesty!"
$ws2.Range("D33").Value = " Ррррааа![K] Мы не смогли пройти\nподземелье! Это издевательство!"
$ws2.Range("E33").Value = " Ññññààà![K] Íú îå òíïãìé ðñïêóé\nðïäèåíåìûå! Üóï éèäåâàóåìûòóâï!"
$ws2.Rows.Item(33).RowHeight = 43.2

# --- Row 34 ---
$ws2.Range("A34").Value = "SCRIPT/D25P11A/um2307.ssb"
$ws2.Range("B34").Value = 253
$ws2.Range("C34").Value = " It must have been a mistake of\nsome sort!"
$ws2.Range("D34").Value = " Должно быть, мы совершили\nкакую-то ошибку!"
$ws2.Range("E34").Value = " Äïìçîï áúóû, íú òïâåñšéìé\nëàëôý-óï ïšéáëô!"
$ws2.Rows.Item(34).RowHeight = 43.2

# --- Row 35: end of this group -> border styles 8/10, no script-name cell ---
$ws1.Range("A9:E9").Copy()
$ws2.Range("A35:E35").PasteSpecial(-4122)
$ws2.Range("B35").Value = 256
$ws2.Range("C35").Value = " We are off once again![K]\nThis time, we won't fail!"
$ws2.Range("D35").Value = " Мы попробуем снова![K] На этот\nраз, у нас всё получится!"
$ws2.Range("E35").Value = " Íú ðïðñïáôåí òîïâà![K] Îà üóïó\nñàè, ô îàò âòæ ðïìôœéóòÿ!"
$ws2.Rows.Item(35).RowHeight = 21.6

# --- Row 36: new script group "SCRIPT/H01P99A/um2401.ssb" ---
$ws2.Range("A36").Value = "SCRIPT/H01P99A/um2401.ssb"
$ws2.Range("B36").Value = 194
$ws2.Range("C36").Value = " You have our support too!"
$ws2.Range("D36").Value = " Мы поддерживаем вас!"
$ws2.Range("E36").Value = " Íú ðïääåñçéâàåí âàò!"
$ws2.Rows.Item(36).RowHeight = 43.2

# --- Row 37 (no A cell) ---
$ws2.Range("B37").Value = 197
$ws2.Range("C37").Value = " O sea! You too!"
$ws2.Range("D37").Value = " О море! И ты тоже!"
$ws2.Range("E37").Value = " Ï íïñå! É óú óïçå!"

# --- Row 38 (no A/B cells) ---
$ws2.Range("C38").Value = " O sea![K] You agree!"
$ws2.Range("D38").Value = " О море![K] Ты согласно!"
$ws2.Range("E38").Value = " Ï íïñå![K] Óú òïãìàòîï!"

# --- Row 39 (no A/B cells) ---
$ws2.Range("C39").Value = " You offer encouragement too!"
$ws2.Range("D39").Value = " Ты тоже их поддерживаешь!"
$ws2.Range("E39").Value = " Óú óïçå éö ðïääåñçéâàåšû!"

# --- Row 40 (no A/B cells) - reuses existing shared strings verbatim ---
$ws2.Range("C40").Value = "[CN]SPLAAAASH…"
$ws2.Range("D40").Value = "[CN]Ш-Ш-Ш-ш-ш-ш..."
$ws2.Range("E40").Value = "[CN]Š-Š-Š-š-š-š..."

# --- Row 41 (no A/B cells) ---
$ws2.Range("C41").Value = "[CN]SPLOOOOSH…"
$ws2.Range("D41").Value = "[CN]Ш-Ш-Ш-Ш-ш-ш..."
$ws2.Range("E41").Value = "[CN]Š-Š-Š-Š-š-š..."

# --- Row 42 (no A/B cells) ---
$ws2.Range("C42").Value = " See?![K] Now do your best!"
$ws2.Range("D42").Value = " Видите?![K] Выложитесь на полную!"
$ws2.Range("E42").Value = " Âéäéóå?![K] Âúìïçéóåòû îà ðïìîôý!"

# --- View: scroll/selection moved to the new bottom of the sheet ---
$ws2.Activate()
$ws2.Range("D43").Select()
$excel.ActiveWindow.ScrollRow = 34

Write-Output "done"
